# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" list in column E (rows 16-28) was re-sorted from
# descending (2102 .. 2002) to ascending (2002 .. 2102) order, and the
# "Valor Mora" figures in column F follow the same reordering (the value
# that used to sit on the first period now sits on the last one, and
# vice-versa).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending order for the "Periodo Mora" column (E16:E28)
$periods = @("2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

# "Valor Mora" column (F16:F28) follows the same reordering: the value that
# used to belong to the first row now belongs to the last row (all the
# values in between were already identical).
$ws.Range("F16").Value = 35112
$ws.Range("F28").Value = 25749
